$wb = $excel.ActiveWorkbook

# --- Metadata sheet (sheet1) ---
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B) to the new timestamp.
$ws.Range("B8").Value = "2024-09-12T14:01:50+00:00"

# Make room for a new "Jurisdiction" property row above the current row 11
# ("Description") by shifting the existing rows 11-21 down into rows
# 12-22 (working bottom-up so values are not overwritten before they are
# copied). Every one of those rows already shares the same cell style, so
# moving values this way (instead of Rows.Insert) keeps formatting intact
# without minting new style records.
for ($r = 21; $r -ge 11; $r--) {
    $dest = $r + 1
    $ws.Range("A$dest").Value = $ws.Range("A$r").Value()
    $ws.Range("B$dest").Value = $ws.Range("B$r").Value()
}

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
